$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price strings (column D),
# matching the workbook's existing text-cell formatting.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '56.217.11'
$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('D3').Value = '2.994.67'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '504.29'
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').Value = '137.92'
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.428'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '7.11'
$ws.Range('E9').Value = '  -3.03%  '
$ws.Range('D10').Value = '0.106'
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('D11').Value = '0.364'
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('D12').Value = '3.501.90'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').Value = '26.10'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '0.0000161'
$ws.Range('E15').Value = '  +1.31%  '
$ws.Range('D16').Value = '56.231.47'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.003.99'
$ws.Range('E17').Value = '  +1.06%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '6.02'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').Value = '12.96'
$ws.Range('E19').Value = '  +3.03%  '
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('D21').Value = '328.35'
$ws.Range('E21').Value = '  +2.93%  '
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').Value = '0.494'
$ws.Range('E23').Value = '  +1.84%  '
$ws.Range('D24').Value = '64.67'
$ws.Range('E24').Value = '  +2.13%  '
$ws.Range('D25').Value = '3.115.70'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  -1.93%  '
$ws.Range('D28').Value = '0.0₃0902'
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('D29').Value = '6.47'
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('D30').Value = '7.01'
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('D31').Value = '1.78'
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').Value = '1.16'
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('D33').Value = '20.22'
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('D34').Value = '152.92'
$ws.Range('E34').Value = '  -1.98%  '
$ws.Range('D35').Value = '4.54'
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('D36').Value = '5.76'
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').Value = '25.54'
$ws.Range('E37').Value = '  +5.73%  '
$ws.Range('D38').Value = '1.26'
$ws.Range('E38').Value = '  +1.26%  '
$ws.Range('D39').Value = '0.0659'
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('D40').Value = '3.038.97'
$ws.Range('E40').Value = '  +1.28%  '
$ws.Range('D41').Value = '36.62'
$ws.Range('E41').Value = '  -2.42%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').Value = '3.80'
$ws.Range('E43').Value = '  +1.79%  '
$ws.Range('E44').Value = '  +2.83%  '
$ws.Range('D45').Value = '2.181.40'
$ws.Range('E45').Value = '  -0.79%  '
$ws.Range('D46').Value = '1.35'
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('D47').Value = '5.90'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').Value = '0.927'
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('D49').Value = '0.0237'
$ws.Range('E49').Value = '  +1.45%  '
$ws.Range('D50').Value = '19.73'
$ws.Range('E50').Value = '  +2.85%  '
$ws.Range('D51').Value = '0.0852'
$ws.Range('E51').Value = '  -2.73%  '
